# Auto-generated COM-interop script applying the scheduled market-data
# refresh described by the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 167006.5
$ws.Range("I55").Value = 250262.25
$ws.Range("J55").Value = 495
$ws.Range("K55").Value = 250262.25
$ws.Range("L55").Value = 495
$ws.Range("M55").Value = -250048.25
$ws.Range("N55").Value = -923

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 24391474
$ws.Range("I137").Value = 33334198
$ws.Range("J137").Value = 2224.182
$ws.Range("K137").Value = 100002594
$ws.Range("L137").Value = 6672.545999999999
$ws.Range("M137").Value = -100000044
$ws.Range("N137").Value = -11772.546

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1238.5333
$ws.Range("I2").Value = 812
$ws.Range("K2").Value = 812
$ws.Range("M2").Value = -699

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2557.1
$ws.Range("I61").Value = 1837.3478
$ws.Range("J61").Value = 4922
$ws.Range("K61").Value = 1837.3478
$ws.Range("L61").Value = 4922
$ws.Range("M61").Value = -1625.3478
$ws.Range("N61").Value = -5346

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 10254.934
$ws.Range("I74").Value = 2502.4
$ws.Range("J74").Value = 25760
$ws.Range("K74").Value = 2502.4
$ws.Range("L74").Value = 25760
$ws.Range("M74").Value = -1628.4
$ws.Range("N74").Value = -27508

$ws = $wb.Worksheets.Item("ARM")
# Row 77
$ws.Range("H77").Value = 10254.934
$ws.Range("I77").Value = 2502.4
$ws.Range("J77").Value = 25760
$ws.Range("K77").Value = 12512
$ws.Range("L77").Value = 128800
$ws.Range("M77").Value = -8144
$ws.Range("N77").Value = -137536

$ws = $wb.Worksheets.Item("ARM")
# Row 116
$ws.Range("H116").Value = 1238.5333
$ws.Range("I116").Value = 812
$ws.Range("K116").Value = 812
$ws.Range("M116").Value = 1482

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2731.261
$ws.Range("I132").Value = 2390.6843
$ws.Range("J132").Value = 4349
$ws.Range("K132").Value = 7172.0529
$ws.Range("L132").Value = 13047
$ws.Range("M132").Value = -4642.0529
$ws.Range("N132").Value = -18107

$ws = $wb.Worksheets.Item("ARM")
# Row 136
$ws.Range("H136").Value = 2557.1
$ws.Range("I136").Value = 1837.3478
$ws.Range("J136").Value = 4922
$ws.Range("K136").Value = 5512.0434
$ws.Range("L136").Value = 14766
$ws.Range("M136").Value = -2962.0434
$ws.Range("N136").Value = -19866

$ws = $wb.Worksheets.Item("ARM")
# Row 139
$ws.Range("H139").Value = 52571.668
$ws.Range("J139").Value = 52571.668
$ws.Range("L139").Value = 52571.668
$ws.Range("N139").Value = -62851.668

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1238.5333
$ws.Range("I3").Value = 812
$ws.Range("K3").Value = 812
$ws.Range("M3").Value = -698

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 19301.334
$ws.Range("I86").Value = 1933.6666
$ws.Range("J86").Value = 36669
$ws.Range("K86").Value = 1933.6666
$ws.Range("L86").Value = 36669
$ws.Range("M86").Value = -810.6666
$ws.Range("N86").Value = -38915

$ws = $wb.Worksheets.Item("BSM")
# Row 89
$ws.Range("H89").Value = 19301.334
$ws.Range("I89").Value = 1933.6666
$ws.Range("J89").Value = 36669
$ws.Range("K89").Value = 9668.333000000001
$ws.Range("L89").Value = 183345
$ws.Range("M89").Value = -4052.333000000001
$ws.Range("N89").Value = -194577

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 886.2143
$ws.Range("I99").Value = 886.2143
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 886.2143
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 611.7857
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4987.467
$ws.Range("I134").Value = 4181.5557
$ws.Range("K134").Value = 12544.6671
$ws.Range("M134").Value = -10009.6671

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 730
$ws.Range("J16").Value = 874.75
$ws.Range("L16").Value = 874.75
$ws.Range("N16").Value = -1448.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1207.7667
$ws.Range("I31").Value = 1133.32
$ws.Range("J31").Value = 1580
$ws.Range("K31").Value = 1133.32
$ws.Range("L31").Value = 1580
$ws.Range("M31").Value = -838.3199999999999
$ws.Range("N31").Value = -2170

$ws = $wb.Worksheets.Item("CRP")
# Row 34
$ws.Range("H34").Value = 1207.7667
$ws.Range("I34").Value = 1133.32
$ws.Range("J34").Value = 1580
$ws.Range("K34").Value = 1133.32
$ws.Range("L34").Value = 1580
$ws.Range("M34").Value = -931.3199999999999
$ws.Range("N34").Value = -1984

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2270.9546
$ws.Range("I58").Value = 1540.8572
$ws.Range("K58").Value = 1540.8572
$ws.Range("M58").Value = -1337.8572

$ws = $wb.Worksheets.Item("CRP")
# Row 113
$ws.Range("H113").Value = 730
$ws.Range("J113").Value = 874.75
$ws.Range("L113").Value = 874.75
$ws.Range("N113").Value = -5214.75

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3173.762
$ws.Range("I132").Value = 2669.111
$ws.Range("J132").Value = 3552.25
$ws.Range("K132").Value = 8007.333
$ws.Range("L132").Value = 10656.75
$ws.Range("M132").Value = -5477.333
$ws.Range("N132").Value = -15716.75

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 2650.5
$ws.Range("I134").Value = 1086.3636
$ws.Range("J134").Value = 6091.6
$ws.Range("K134").Value = 3259.0908
$ws.Range("L134").Value = 18274.8
$ws.Range("M134").Value = -724.0907999999999
$ws.Range("N134").Value = -23344.8

$ws = $wb.Worksheets.Item("CRP")
# Row 136
$ws.Range("H136").Value = 2270.9546
$ws.Range("I136").Value = 1540.8572
$ws.Range("K136").Value = 4622.571599999999
$ws.Range("M136").Value = -2072.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 136
$ws.Range("H136").Value = 4100.606
$ws.Range("I136").Value = 2086
$ws.Range("J136").Value = 4460.357
$ws.Range("K136").Value = 6258
$ws.Range("L136").Value = 13381.071
$ws.Range("M136").Value = -1158
$ws.Range("N136").Value = -23581.071

$ws = $wb.Worksheets.Item("GSM")
# Row 121
$ws.Range("H121").Value = 25982
$ws.Range("J121").Value = 25982
$ws.Range("L121").Value = 25982
$ws.Range("N121").Value = -29476

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2583.5
$ws.Range("I132").Value = 1576.5714
$ws.Range("J132").Value = 4933
$ws.Range("K132").Value = 4729.7142
$ws.Range("L132").Value = 14799
$ws.Range("M132").Value = -2199.7142
$ws.Range("N132").Value = -19859

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 7805.1333
$ws.Range("I22").Value = 917.1667
$ws.Range("J22").Value = 12397.111
$ws.Range("K22").Value = 917.1667
$ws.Range("L22").Value = 12397.111
$ws.Range("M22").Value = -622.1667
$ws.Range("N22").Value = -12987.111

$ws = $wb.Worksheets.Item("LTW")
# Row 27
$ws.Range("H27").Value = 7805.1333
$ws.Range("I27").Value = 917.1667
$ws.Range("J27").Value = 12397.111
$ws.Range("K27").Value = 917.1667
$ws.Range("L27").Value = 12397.111
$ws.Range("M27").Value = -810.1667
$ws.Range("N27").Value = -12611.111

$ws = $wb.Worksheets.Item("LTW")
# Row 29
$ws.Range("H29").Value = 16129
$ws.Range("I29").Value = 18838.666
$ws.Range("K29").Value = 18838.666
$ws.Range("M29").Value = -18543.666

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2390.0476
$ws.Range("I100").Value = 1719.1
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1719.1
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1178.1
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3397.6453
$ws.Range("I132").Value = 2225.7
$ws.Range("J132").Value = 5528.4546
$ws.Range("K132").Value = 6677.099999999999
$ws.Range("L132").Value = 16585.3638
$ws.Range("M132").Value = -4147.099999999999
$ws.Range("N132").Value = -21645.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 5141.8
$ws.Range("I136").Value = 3255.6667
$ws.Range("J136").Value = 7971
$ws.Range("K136").Value = 9767.000100000001
$ws.Range("L136").Value = 23913
$ws.Range("M136").Value = -7217.000100000001
$ws.Range("N136").Value = -29013

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 10419057
$ws.Range("I132").Value = 12822312
$ws.Range("J132").Value = 4947.6665
$ws.Range("K132").Value = 38466936
$ws.Range("L132").Value = 14842.9995
$ws.Range("M132").Value = -38464406
$ws.Range("N132").Value = -19902.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 17598636
$ws.Range("I136").Value = 37149228
$ws.Range("K136").Value = 111447684
$ws.Range("M136").Value = -111445134
$ws.Range("N136").Value = -14428.7
